$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 467, shifting rows 467:515 down to 468:516
$ws.Rows("467:467").Insert()

# Populate the new row 467 with the new record
$ws.Range("A467").Value = 4
$ws.Range("B467").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C467").Value = "Los Lagos"
$ws.Range("D467").Value = 45212
$ws.Range("D467").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E467").Value = 10
$ws.Range("F467").Value = 100112043
$ws.Range("G467").Value = "Pepino ensalada"
$ws.Range("H467").Value = "Sin especificar"
$ws.Range("I467").Value = "Primera"
$ws.Range("J467").Value = 400
$ws.Range("K467").Value = 17000
$ws.Range("L467").Value = 17000
$ws.Range("M467").Value = 17000
$ws.Range("N467").Value = "$/caja 60 unidades"
$ws.Range("O467").Value = "Región de Arica y Parinacota"
$ws.Range("P467").Value = 283
$ws.Range("Q467").Value = 60
$ws.Range("R467").Value = "Hortaliza"
